# Insert a new data row after the existing header+96 data rows (i.e. at
# row 98), shifting all the rows that were 98:221 down to 99:222, then
# fill the newly-inserted row 98 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 98 (and everything below it) down by one row.
$ws.Rows(98).Insert()

# Populate the freshly inserted row with the new record.
$ws.Cells.Item(98, 1).Value  = 4
$ws.Cells.Item(98, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(98, 3).Value  = "Los Lagos"
$ws.Cells.Item(98, 4).Value  = 44579
$ws.Cells.Item(98, 5).Value  = 10
$ws.Cells.Item(98, 6).Value  = 100112037
$ws.Cells.Item(98, 7).Value  = "Cebollín"
$ws.Cells.Item(98, 8).Value  = "Sin especificar"
$ws.Cells.Item(98, 9).Value  = "Primera"
$ws.Cells.Item(98, 10).Value = 160
$ws.Cells.Item(98, 11).Value = 6000
$ws.Cells.Item(98, 12).Value = 6500
$ws.Cells.Item(98, 13).Value = 6250
$ws.Cells.Item(98, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(98, 15).Value = "Región Metropolitana"
$ws.Cells.Item(98, 16).Value = 174
$ws.Cells.Item(98, 17).Value = 36
$ws.Cells.Item(98, 18).Value = "Hortaliza"
